$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.5052451069971653
$ws.Cells.Item(2, 3).Value = 0.03838363311821524
$ws.Cells.Item(2, 4).Value = 0.1762427502561366
$ws.Cells.Item(2, 6).Value = 2.020745369541231
$ws.Cells.Item(2, 7).Value = 1.292789114319987
$ws.Cells.Item(2, 8).Value = 1.225385994647766
$ws.Cells.Item(2, 10).Value = 0.2849123555693325
$ws.Cells.Item(2, 11).Value = 0.4663479685445395
$ws.Cells.Item(2, 13).Value = 0.3243925944714405
$ws.Cells.Item(2, 14).Value = 2.261013135025742
$ws.Cells.Item(3, 2).Value = 0.470537953993329
$ws.Cells.Item(3, 3).Value = 0.0337565837968441
$ws.Cells.Item(3, 4).Value = 0.173065920998198
$ws.Cells.Item(3, 6).Value = 2.018328926407705
$ws.Cells.Item(3, 7).Value = 1.290651734703999
$ws.Cells.Item(3, 8).Value = 1.2291682874461
$ws.Cells.Item(3, 10).Value = 0.2830260406512721
$ws.Cells.Item(3, 11).Value = 0.4296111438855803
$ws.Cells.Item(3, 13).Value = 0.3120357957168736
$ws.Cells.Item(3, 14).Value = 2.280554701740044
$ws.Cells.Item(4, 2).Value = 0.4494599612741581
$ws.Cells.Item(4, 3).Value = 0.03091873727477434
$ws.Cells.Item(4, 4).Value = 0.1711892472870602
$ws.Cells.Item(4, 6).Value = 2.017806375024534
$ws.Cells.Item(4, 7).Value = 1.290027140437857
$ws.Cells.Item(4, 8).Value = 1.231978801297117
$ws.Cells.Item(4, 10).Value = 0.2820036415501974
$ws.Cells.Item(4, 11).Value = 0.4072482278075142
$ws.Cells.Item(4, 13).Value = 0.3046148592141407
$ws.Cells.Item(4, 14).Value = 2.293230099939734
$ws.Cells.Item(5, 2).Value = 0.4409291969961942
$ws.Cells.Item(5, 3).Value = 0.02976310740700683
$ws.Cells.Item(5, 4).Value = 0.1704431459982345
$ws.Cells.Item(5, 6).Value = 2.017835180667376
$ws.Cells.Item(5, 7).Value = 1.289945464361963
$ws.Cells.Item(5, 8).Value = 1.233246930721876
$ws.Cells.Item(5, 10).Value = 0.2816211724433373
$ws.Cells.Item(5, 11).Value = 0.3981840955030975
$ws.Cells.Item(5, 13).Value = 0.3016326633247886
$ws.Cells.Item(5, 14).Value = 2.298565568594992
$ws.Cells.Item(6, 2).Value = 0.4395162226094556
$ws.Cells.Item(6, 3).Value = 0.02957126568728086
$ws.Cells.Item(6, 4).Value = 0.170320385516618
$ws.Cells.Item(6, 6).Value = 2.01785456671648
$ws.Cells.Item(6, 7).Value = 1.289942339217717
$ws.Cells.Item(6, 8).Value = 1.233464922824666
$ws.Cells.Item(6, 10).Value = 0.281559727853022
$ws.Cells.Item(6, 11).Value = 0.39668196620255
$ws.Cells.Item(6, 13).Value = 0.3011400055459461
$ws.Cells.Item(6, 14).Value = 2.299461793920514
$ws.Cells.Item(7, 2).Value = 0.4493446744216101
$ws.Cells.Item(7, 3).Value = 0.03090314872007127
$ws.Cells.Item(7, 4).Value = 0.1711791094594872
$ws.Cells.Item(7, 6).Value = 2.017805784597087
$ws.Cells.Item(7, 7).Value = 1.290025339183629
$ws.Cells.Item(7, 8).Value = 1.231995406344907
$ws.Cells.Item(7, 10).Value = 0.2819983450742853
$ws.Cells.Item(7, 11).Value = 0.4071257872778915
$ws.Cells.Item(7, 13).Value = 0.3045744705000217
$ws.Cells.Item(7, 14).Value = 2.293301367106356
$ws.Cells.Item(8, 2).Value = 0.4932300057374448
$ws.Cells.Item(8, 3).Value = 0.0367875744708499
$ws.Cells.Item(8, 4).Value = 0.175132081194306
$ws.Cells.Item(8, 6).Value = 2.019712721741513
$ws.Cells.Item(8, 7).Value = 1.291909339715062
$ws.Cells.Item(8, 8).Value = 1.226588848177514
$ws.Cells.Item(8, 10).Value = 0.2842337767326342
$ws.Cells.Item(8, 11).Value = 0.4536410551461074
$ws.Cells.Item(8, 13).Value = 0.3200975372610486
$ws.Cells.Item(8, 14).Value = 2.267610572249353
$ws.Cells.Item(9, 2).Value = 0.5811257300912871
$ws.Cells.Item(9, 3).Value = 0.04835204293719642
$ws.Cells.Item(9, 4).Value = 0.1834676823901589
$ws.Cells.Item(9, 6).Value = 2.031078689461125
$ws.Cells.Item(9, 7).Value = 1.301067119789622
$ws.Cells.Item(9, 8).Value = 1.219857824867162
$ws.Cells.Item(9, 10).Value = 0.2896948257275298
$ws.Cells.Item(9, 11).Value = 0.5463881163709345
$ws.Cells.Item(9, 13).Value = 0.3518541012631999
$ws.Cells.Item(9, 14).Value = 2.222602485651691
$ws.Cells.Item(10, 2).Value = 0.646820505845767
$ws.Cells.Item(10, 3).Value = 0.05686450344992977
$ws.Cells.Item(10, 4).Value = 0.1899450774377982
$ws.Cells.Item(10, 6).Value = 2.044082842158929
$ws.Cells.Item(10, 7).Value = 1.311137225141266
$ws.Cells.Item(10, 8).Value = 1.2172707484916
$ws.Cells.Item(10, 10).Value = 0.2943644628007576
$ws.Cells.Item(10, 11).Value = 0.6154623539423767
$ws.Cells.Item(10, 13).Value = 0.3759870153358378
$ws.Cells.Item(10, 14).Value = 2.192810674596764
$ws.Cells.Item(11, 2).Value = 0.6769495222191892
$ws.Cells.Item(11, 3).Value = 0.06074073663150159
$ws.Cells.Item(11, 4).Value = 0.1929679960919088
$ws.Cells.Item(11, 6).Value = 2.051010471387826
$ws.Cells.Item(11, 7).Value = 1.316446678256384
$ws.Cells.Item(11, 8).Value = 1.216605611551174
$ws.Cells.Item(11, 10).Value = 0.2966317407612848
$ws.Cells.Item(11, 11).Value = 0.6470890231794613
$ws.Cells.Item(11, 13).Value = 0.3871397210936962
$ws.Cells.Item(11, 14).Value = 2.179969128978861
$ws.Cells.Item(12, 2).Value = 0.688393542662709
$ws.Cells.Item(12, 3).Value = 0.06220911955298902
$ws.Cells.Item(12, 4).Value = 0.1941236104323991
$ws.Cells.Item(12, 6).Value = 2.053779333337829
$ws.Cells.Item(12, 7).Value = 1.318562153026704
$ws.Cells.Item(12, 8).Value = 1.216427291577844
$ws.Cells.Item(12, 10).Value = 0.2975108671480911
$ws.Cells.Item(12, 11).Value = 0.6590944849716038
$ws.Cells.Item(12, 13).Value = 0.3913879935505307
$ws.Cells.Item(12, 14).Value = 2.175208659676883
$ws.Cells.Item(13, 2).Value = 0.6859273242757524
$ws.Cells.Item(13, 3).Value = 0.06189285341666562
$ws.Cells.Item(13, 4).Value = 0.1938742445744026
$ws.Cells.Item(13, 6).Value = 2.053176537383393
$ws.Cells.Item(13, 7).Value = 1.318101879968239
$ws.Cells.Item(13, 8).Value = 1.216462425218268
$ws.Cells.Item(13, 10).Value = 0.2973206174281557
$ws.Cells.Item(13, 11).Value = 0.6565076027870589
$ws.Cells.Item(13, 13).Value = 0.3904719426279044
$ws.Cells.Item(13, 14).Value = 2.176229358149897
$ws.Cells.Item(14, 2).Value = 0.6778903325694614
$ws.Cells.Item(14, 3).Value = 0.06086153071510125
$ws.Cells.Item(14, 4).Value = 0.1930628511370571
$ws.Cells.Item(14, 6).Value = 2.051235350911583
$ws.Cells.Item(14, 7).Value = 1.316618616473932
$ws.Cells.Item(14, 8).Value = 1.216589467140906
$ws.Cells.Item(14, 10).Value = 0.2967036550792841
$ws.Cells.Item(14, 11).Value = 0.6480761375204054
$ws.Cells.Item(14, 13).Value = 0.3874887289287727
$ws.Cells.Item(14, 14).Value = 2.179575430458939
$ws.Cells.Item(15, 2).Value = 0.6729719685242799
$ws.Cells.Item(15, 3).Value = 0.06022988521465322
$ws.Cells.Item(15, 4).Value = 0.1925672666230582
$ws.Cells.Item(15, 6).Value = 2.050065269006453
$ws.Cells.Item(15, 7).Value = 1.315723740045428
$ws.Cells.Item(15, 8).Value = 1.216676861757108
$ws.Cells.Item(15, 10).Value = 0.2963284246848872
$ws.Cells.Item(15, 11).Value = 0.6429154069375898
$ws.Cells.Item(15, 13).Value = 0.3856646748638539
$ws.Cells.Item(15, 14).Value = 2.181638328332646
$ws.Cells.Item(16, 2).Value = 0.6448563915387524
$ws.Cells.Item(16, 3).Value = 0.05661125876376616
$ws.Cells.Item(16, 4).Value = 0.189749052006448
$ws.Cells.Item(16, 6).Value = 2.043650468643818
$ws.Cells.Item(16, 7).Value = 1.310804911445032
$ws.Cells.Item(16, 8).Value = 1.217324511259477
$ws.Cells.Item(16, 10).Value = 0.2942191681673734
$ws.Cells.Item(16, 11).Value = 0.6133995676168524
$ws.Cells.Item(16, 13).Value = 0.3752616621042577
$ws.Cells.Item(16, 14).Value = 2.193664214995188
$ws.Cells.Item(17, 2).Value = 0.6276707175340164
$ws.Cells.Item(17, 3).Value = 0.05439232665409577
$ws.Cells.Item(17, 4).Value = 0.1880396631165269
$ws.Cells.Item(17, 6).Value = 2.039974385813949
$ws.Cells.Item(17, 7).Value = 1.307974058796034
$ws.Cells.Item(17, 8).Value = 1.217852866657651
$ws.Cells.Item(17, 10).Value = 0.2929618347987315
$ws.Cells.Item(17, 11).Value = 0.5953447280599278
$ws.Cells.Item(17, 13).Value = 0.3689243708281822
$ws.Cells.Item(17, 14).Value = 2.201223868746666
$ws.Cells.Item(18, 2).Value = 0.6178089687129784
$ws.Cells.Item(18, 3).Value = 0.05311642120685178
$ws.Cells.Item(18, 4).Value = 0.1870636557713539
$ws.Cells.Item(18, 6).Value = 2.037955252059604
$ws.Cells.Item(18, 7).Value = 1.306414389696116
$ws.Cells.Item(18, 8).Value = 1.218204931778203
$ws.Cells.Item(18, 10).Value = 0.2922521138376482
$ws.Cells.Item(18, 11).Value = 0.5849793220824893
$ws.Cells.Item(18, 13).Value = 0.3652957619886124
$ws.Cells.Item(18, 14).Value = 2.205638887439164
$ws.Cells.Item(19, 2).Value = 0.6144739085003152
$ws.Cells.Item(19, 3).Value = 0.05268448492958555
$ws.Cells.Item(19, 4).Value = 0.1867344332446947
$ws.Cells.Item(19, 6).Value = 2.037287967744092
$ws.Cells.Item(19, 7).Value = 1.305898083930927
$ws.Cells.Item(19, 8).Value = 1.218332409130028
$ws.Cells.Item(19, 10).Value = 0.2920141272989554
$ws.Cells.Item(19, 11).Value = 0.5814730899558356
$ws.Cells.Item(19, 13).Value = 0.3640700032426878
$ws.Cells.Item(19, 14).Value = 2.207145224481444
$ws.Cells.Item(20, 2).Value = 0.6294977839406783
$ws.Cells.Item(20, 3).Value = 0.05462849806504266
$ws.Cells.Item(20, 4).Value = 0.1882208870943742
$ws.Cells.Item(20, 6).Value = 2.040355852908974
$ws.Cells.Item(20, 7).Value = 1.308268310912908
$ws.Cells.Item(20, 8).Value = 1.217791637253768
$ws.Cells.Item(20, 10).Value = 0.2930942866962027
$ws.Cells.Item(20, 11).Value = 0.5972647051623312
$ws.Cells.Item(20, 13).Value = 0.3695972863058543
$ws.Cells.Item(20, 14).Value = 2.200412205342005
$ws.Cells.Item(21, 2).Value = 0.6802500492735248
$ws.Cells.Item(21, 3).Value = 0.06116444075665584
$ws.Cells.Item(21, 4).Value = 0.1933008818061239
$ws.Cells.Item(21, 6).Value = 2.051801575092171
$ws.Cells.Item(21, 7).Value = 1.317051438829779
$ws.Cells.Item(21, 8).Value = 1.216550155909033
$ws.Cells.Item(21, 10).Value = 0.2968843140412076
$ws.Cells.Item(21, 11).Value = 0.6505518748721215
$ws.Cells.Item(21, 13).Value = 0.3883642939109038
$ws.Cells.Item(21, 14).Value = 2.178589829324046
$ws.Cells.Item(22, 2).Value = 0.7136223873370966
$ws.Cells.Item(22, 3).Value = 0.06543918148186378
$ws.Cells.Item(22, 4).Value = 0.1966844551555482
$ws.Cells.Item(22, 6).Value = 2.060130213240171
$ws.Cells.Item(22, 7).Value = 1.323403224837676
$ws.Cells.Item(22, 8).Value = 1.216167483981451
$ws.Cells.Item(22, 10).Value = 0.2994811354283087
$ws.Cells.Item(22, 11).Value = 0.6855477784466188
$ws.Cells.Item(22, 13).Value = 0.4007751974880307
$ws.Cells.Item(22, 14).Value = 2.16492424953028
$ws.Cells.Item(23, 2).Value = 0.6957924821772963
$ws.Cells.Item(23, 3).Value = 0.06315739236572426
$ws.Cells.Item(23, 4).Value = 0.1948727922283098
$ws.Cells.Item(23, 6).Value = 2.055607452391115
$ws.Cells.Item(23, 7).Value = 1.319957157610702
$ws.Cells.Item(23, 8).Value = 1.216332507376706
$ws.Cells.Item(23, 10).Value = 0.2980842031058017
$ws.Cells.Item(23, 11).Value = 0.6668543719717661
$ws.Cells.Item(23, 13).Value = 0.3941379808739427
$ws.Cells.Item(23, 14).Value = 2.172163197149665
$ws.Cells.Item(24, 2).Value = 0.6286717094958192
$ws.Cells.Item(24, 3).Value = 0.05452172562894475
$ws.Cells.Item(24, 4).Value = 0.1881389347254014
$ws.Cells.Item(24, 6).Value = 2.040183097905299
$ws.Cells.Item(24, 7).Value = 1.308135068263809
$ws.Cells.Item(24, 8).Value = 1.21781916858663
$ws.Cells.Item(24, 10).Value = 0.2930343642734954
$ws.Cells.Item(24, 11).Value = 0.596396638137179
$ws.Cells.Item(24, 13).Value = 0.3692930151595277
$ws.Cells.Item(24, 14).Value = 2.200778943828837
$ws.Cells.Item(25, 2).Value = 0.5571511750769105
$ws.Cells.Item(25, 3).Value = 0.04522078305845412
$ws.Cells.Item(25, 4).Value = 0.1811504568678259
$ws.Cells.Item(25, 6).Value = 2.027186992511133
$ws.Cells.Item(25, 7).Value = 1.29800365049438
$ws.Cells.Item(25, 8).Value = 1.221264494628869
$ws.Cells.Item(25, 10).Value = 0.2881020648318469
$ws.Cells.Item(25, 11).Value = 0.5211337253785189
$ws.Cells.Item(25, 13).Value = 0.3431223220390294
$ws.Cells.Item(25, 14).Value = 2.234203217887298
